$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44215
$ws.Cells.Item(2, 10).Value = 16000

# Row 3
$ws.Cells.Item(3, 4).Value = 44231
$ws.Cells.Item(3, 10).Value = 12000

# Row 4
$ws.Cells.Item(4, 4).Value = 44160

# Row 5
$ws.Cells.Item(5, 4).Value = 44214
$ws.Cells.Item(5, 10).Value = 7000

# Row 6
$ws.Cells.Item(6, 4).Value = 44167
$ws.Cells.Item(6, 10).Value = 7000

# Row 7
$ws.Cells.Item(7, 4).Value = 44845

# Row 8
$ws.Cells.Item(8, 4).Value = 44880
$ws.Cells.Item(8, 10).Value = 7900

# Row 9
$ws.Cells.Item(9, 4).Value = 44188
$ws.Cells.Item(9, 10).Value = 12000

# Row 10
$ws.Cells.Item(10, 4).Value = 44166
$ws.Cells.Item(10, 10).Value = 7000

# Row 11
$ws.Cells.Item(11, 4).Value = 44187
$ws.Cells.Item(11, 10).Value = 12000

# Row 12
$ws.Cells.Item(12, 4).Value = 44210
$ws.Cells.Item(12, 10).Value = 8800
$ws.Cells.Item(12, 11).Value = 2500
$ws.Cells.Item(12, 13).Value = 2750
$ws.Cells.Item(12, 16).Value = 28

# Row 13
$ws.Cells.Item(13, 4).Value = 44209
$ws.Cells.Item(13, 11).Value = 2500
$ws.Cells.Item(13, 13).Value = 2750
$ws.Cells.Item(13, 16).Value = 28

# Row 14
$ws.Cells.Item(14, 4).Value = 44846
$ws.Cells.Item(14, 10).Value = 7900
$ws.Cells.Item(14, 11).Value = 3000
$ws.Cells.Item(14, 13).Value = 3000
$ws.Cells.Item(14, 16).Value = 30

# Row 15
$ws.Cells.Item(15, 4).Value = 44189
$ws.Cells.Item(15, 10).Value = 16000

# Row 16
$ws.Cells.Item(16, 4).Value = 44162
$ws.Cells.Item(16, 10).Value = 7000
$ws.Cells.Item(16, 11).Value = 3000
$ws.Cells.Item(16, 12).Value = 3000
$ws.Cells.Item(16, 13).Value = 3000
$ws.Cells.Item(16, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(16, 16).Value = 30

# Row 17
$ws.Cells.Item(17, 4).Value = 44229
$ws.Cells.Item(17, 10).Value = 16000

# Row 18
$ws.Cells.Item(18, 4).Value = 44859
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 7900
$ws.Cells.Item(18, 11).Value = 3000
$ws.Cells.Item(18, 12).Value = 3000
$ws.Cells.Item(18, 13).Value = 3000
$ws.Cells.Item(18, 16).Value = 30

# Row 19
$ws.Cells.Item(19, 4).Value = 44874
$ws.Cells.Item(19, 10).Value = 7900

# Row 20
$ws.Cells.Item(20, 4).Value = 44181
$ws.Cells.Item(20, 10).Value = 12000

# Row 21
$ws.Cells.Item(21, 4).Value = 44876
$ws.Cells.Item(21, 10).Value = 7900

# Row 22
$ws.Cells.Item(22, 4).Value = 44855

# Row 23
$ws.Cells.Item(23, 4).Value = 44602

# Row 24
$ws.Cells.Item(24, 4).Value = 44602
$ws.Cells.Item(24, 9).Value = "Segunda"
$ws.Cells.Item(24, 10).Value = 6000
$ws.Cells.Item(24, 12).Value = 2500
$ws.Cells.Item(24, 13).Value = 2500
$ws.Cells.Item(24, 16).Value = 25

# Row 25
$ws.Cells.Item(25, 4).Value = 44186
$ws.Cells.Item(25, 10).Value = 10000

# Row 26
$ws.Cells.Item(26, 4).Value = 44245
$ws.Cells.Item(26, 10).Value = 9000
$ws.Cells.Item(26, 15).Value = "Región Metropolitana"

# Row 27
$ws.Cells.Item(27, 4).Value = 44245
$ws.Cells.Item(27, 9).Value = "Segunda"
$ws.Cells.Item(27, 10).Value = 5000
$ws.Cells.Item(27, 11).Value = 2500
$ws.Cells.Item(27, 12).Value = 2500
$ws.Cells.Item(27, 13).Value = 2500
$ws.Cells.Item(27, 15).Value = "Región Metropolitana"
$ws.Cells.Item(27, 16).Value = 25

# Row 28
$ws.Cells.Item(28, 4).Value = 44204
$ws.Cells.Item(28, 10).Value = 7000

# Row 29
$ws.Cells.Item(29, 4).Value = 44873
$ws.Cells.Item(29, 10).Value = 7900

# Row 30
$ws.Cells.Item(30, 4).Value = 44860
$ws.Cells.Item(30, 10).Value = 7900
$ws.Cells.Item(30, 15).Value = "Provincia de Chacabuco"

# Row 31
$ws.Cells.Item(31, 4).Value = 44883
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 9700
$ws.Cells.Item(31, 11).Value = 3000
$ws.Cells.Item(31, 12).Value = 3000
$ws.Cells.Item(31, 13).Value = 3000
$ws.Cells.Item(31, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(31, 16).Value = 30

# Row 32
$ws.Cells.Item(32, 4).Value = 44230
$ws.Cells.Item(32, 10).Value = 16000

# Row 33
$ws.Cells.Item(33, 4).Value = 44882
$ws.Cells.Item(33, 10).Value = 7900

# Row 34
$ws.Cells.Item(34, 4).Value = 44159
$ws.Cells.Item(34, 10).Value = 7000

# Row 35
$ws.Cells.Item(35, 4).Value = 44875

# Row 36
$ws.Cells.Item(36, 4).Value = 44600
$ws.Cells.Item(36, 10).Value = 1300
$ws.Cells.Item(36, 11).Value = 3500
$ws.Cells.Item(36, 12).Value = 4000
$ws.Cells.Item(36, 13).Value = 3808
$ws.Cells.Item(36, 15).Value = "Región Metropolitana"
$ws.Cells.Item(36, 16).Value = 38

# Row 37
$ws.Cells.Item(37, 4).Value = 44847

# Row 38
$ws.Cells.Item(38, 4).Value = 44168
$ws.Cells.Item(38, 10).Value = 7000

# Row 39
$ws.Cells.Item(39, 4).Value = 44881
$ws.Cells.Item(39, 10).Value = 7900

# Row 40
$ws.Cells.Item(40, 4).Value = 44161

# Row 41
$ws.Cells.Item(41, 4).Value = 44232
$ws.Cells.Item(41, 10).Value = 16000
